$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column A stays text (prevents Excel from auto-converting
# dash-separated dates like "01-08-2022" into date serials)
$ws.Range("A3:A21").NumberFormat = "@"

# Row 3: date format change + D/G updated to 1 (E/F/H unchanged)
$ws.Range("A3").Value = "28-07-2022"
$ws.Range("D3").Value = 1
$ws.Range("G3").Value = 1

# Row 4: date format change + D/E set to 1, H set to 0
$ws.Range("A4").Value = "01-08-2022"
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 1
$ws.Range("H4").Value = 0

# Row 5
$ws.Range("A5").Value = "04-08-2022"
$ws.Range("D5").Value = 1
$ws.Range("E5").Value = 1
$ws.Range("H5").Value = 0

# Row 6
$ws.Range("A6").Value = "08-08-2022"
$ws.Range("D6").Value = 1
$ws.Range("E6").Value = 1
$ws.Range("H6").Value = 0

# Row 7
$ws.Range("A7").Value = "11-08-2022"
$ws.Range("D7").Value = 1
$ws.Range("E7").Value = 1
$ws.Range("H7").Value = 0

# Row 8: only date format changes
$ws.Range("A8").Value = "15-08-2022"

# Row 9: only date format changes
$ws.Range("A9").Value = "18-08-2022"

# Row 10: only date format changes
$ws.Range("A10").Value = "22-08-2022"

# Row 11
$ws.Range("A11").Value = "25-08-2022"
$ws.Range("D11").Value = 1
$ws.Range("E11").Value = 1
$ws.Range("H11").Value = 0

# Row 12: only date format changes
$ws.Range("A12").Value = "29-08-2022"

# Row 13: only date format changes
$ws.Range("A13").Value = "01-09-2022"

# Row 14
$ws.Range("A14").Value = "05-09-2022"
$ws.Range("D14").Value = 1
$ws.Range("E14").Value = 1
$ws.Range("H14").Value = 0

# Row 15: only date format changes
$ws.Range("A15").Value = "08-09-2022"

# Row 16: only date format changes
$ws.Range("A16").Value = "12-09-2022"

# Row 17: only date format changes
$ws.Range("A17").Value = "15-09-2022"

# Row 18: only date format changes
$ws.Range("A18").Value = "19-09-2022"

# Row 19: only date format changes
$ws.Range("A19").Value = "22-09-2022"

# Row 20
$ws.Range("A20").Value = "26-09-2022"
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = 1
$ws.Range("H20").Value = 0

# Row 21: only date format changes
$ws.Range("A21").Value = "29-09-2022"
